$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking Price values remain text by forcing Text format before assignment
$ws.Range('D2').Value = '26.031.98'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').Value = '1.637.00'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('E4').Value = '  +0.53%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.85'
$ws.Range('E5').Value = '  -0.30%  '
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('E8').Value = '  -1.62%  '
$ws.Range('E9').Value = '  -1.39%  '
$ws.Range('E10').Value = '  -4.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0795'
$ws.Range('E11').Value = '  +0.31%  '
$ws.Range('D12').Value = '1.702.16'
$ws.Range('E12').Value = '  +3.71%  '
$ws.Range('E13').Value = '  -1.57%  '
$ws.Range('E14').Value = '  -1.97%  '
$ws.Range('D15').Value = '0.0₃0749'
$ws.Range('E15').Value = '  -1.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '62.26'
$ws.Range('E16').Value = '  -0.91%  '
$ws.Range('D17').Value = '26.052.85'
$ws.Range('E17').Value = '  +0.33%  '
$ws.Range('E18').Value = '  +0.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '191.35'
$ws.Range('E19').Value = '  -0.82%  '
$ws.Range('E20').Value = '  -1.82%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.64'
$ws.Range('E21').Value = '  -2.76%  '
$ws.Range('E22').Value = '  -1.60%  '
$ws.Range('E23').Value = '  +0.94%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '143.90'
$ws.Range('E24').Value = '  +0.44%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.78'
$ws.Range('E25').Value = '  -0.35%  '
$ws.Range('E26').Value = '  +0.67%  '
$ws.Range('E27').Value = '  -1.51%  '
$ws.Range('E28').Value = '  -1.63%  '
$ws.Range('E29').Value = '  -0.35%  '
$ws.Range('E30').Value = '  -2.95%  '
$ws.Range('E31').Value = '  -2.21%  '
$ws.Range('E32').Value = '  -3.11%  '
$ws.Range('E35').Value = '  -2.25%  '
$ws.Range('D36').Value = '1.131.43'
$ws.Range('E36').Value = '  -0.21%  '
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('E38').Value = '  -2.56%  '
$ws.Range('E39').Value = '  -0.81%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '98.95'
$ws.Range('E40').Value = '  -0.33%  '
$ws.Range('D43').Value = '0.0₆0113'
$ws.Range('E43').Value = '  -1.11%  '
$ws.Range('E44').Value = '  -1.72%  '
$ws.Range('E45').Value = '  -0.60%  '
$ws.Range('E46').Value = '  +1.38%  '
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.60'
$ws.Range('E49').Value = '  +0.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0930'
$ws.Range('E50').Value = '  -2.92%  '
$ws.Range('E51').Value = '  -0.04%  '

Write-Host "Updated cryptos list"
